$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.717.24'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '1.598.70'
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("E6").Value = '  -0.67%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.0619'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.54%  '
$ws.Range("E9").Value = '  +1.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.52'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.63%  '
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("D12").Value = '1.821.93'
$ws.Range("E12").Value = '  +0.18%  '
$ws.Range("D13").Value = '1.605.62'
$ws.Range("E13").Value = '  +1.47%  '
$ws.Range("E14").Value = '  +0.61%  '
$ws.Range("E15").Value = '  +0.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.41'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.49%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '26.671.53'
$ws.Range("E17").Value = '  +0.26%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '0.0₃0765'
$ws.Range("E18").Value = '  +5.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '209.86'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.18%  '
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.15'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.62%  '
$ws.Range("E22").Value = '  +0.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.30'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.88%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.10'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.01'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.31%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("E28").Value = '  +0.16%  '
$ws.Range("E29").Value = '  +0.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0522'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.85%  '
$ws.Range("E31").Value = '  -0.08%  '
$ws.Range("E32").Value = '  +0.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.96'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.07%  '
$ws.Range("D34").Value = '1.286.58'
$ws.Range("E34").Value = '  +0.57%  '
$ws.Range("E35").Value = '  -5.26%  '
$ws.Range("E36").Value = '  +0.48%  '
$ws.Range("E37").Value = '  +0.34%  '
$ws.Range("E38").Value = '  -0.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.08'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +17.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.828'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.26%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.45'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.45%  '
$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.19'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.56%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.783'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.26'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("D45").Value = '1.735.03'
$ws.Range("E45").Value = '  +0.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.33'
$ws.Range("D46").Style = "Normal"
$ws.Range("E47").Value = '  -1.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.100'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("E49").Value = '  +0.60%  '
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("E51").Value = '  -1.65%  '
